$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values that look numeric must be forced to text
# to match the source data which stores them as strings (e.g. "1.00", "18.00").
$dCells = @("D2","D3","D4","D5","D6","D8","D12","D13","D14","D15","D16","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "58.929.79"
$ws.Range("D3").Value = "2.604.64"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "553.17"
$ws.Range("D6").Value = "144.67"
$ws.Range("D8").Value = "0.598"
$ws.Range("D12").Value = "0.336"
$ws.Range("D13").Value = "3.056.05"
$ws.Range("D14").Value = "58.870.82"
$ws.Range("D15").Value = "20.91"
$ws.Range("D16").Value = "2.605.42"
$ws.Range("D18").Value = "4.47"
$ws.Range("D19").Value = "339.05"
$ws.Range("D20").Value = "10.13"
$ws.Range("D21").Value = "6.18"
$ws.Range("D22").Value = "0.998"
$ws.Range("D23").Value = "66.68"
$ws.Range("D24").Value = "0.426"
$ws.Range("D25").Value = "0.993"
$ws.Range("D26").Value = "0.159"
$ws.Range("D27").Value = "7.15"
$ws.Range("D28").Value = "0.0₃0763"
$ws.Range("D31").Value = "5.97"
$ws.Range("D32").Value = "154.47"
$ws.Range("D33").Value = "18.93"
$ws.Range("D34").Value = "3.94"
$ws.Range("D35").Value = "0.879"
$ws.Range("D36").Value = "1.12"
$ws.Range("D37").Value = "37.25"
$ws.Range("D38").Value = "1.46"
$ws.Range("D39").Value = "0.826"
$ws.Range("D40").Value = "3.62"
$ws.Range("D41").Value = "283.69"
$ws.Range("D42").Value = "0.999"
$ws.Range("D43").Value = "0.600"
$ws.Range("D44").Value = "0.0953"
$ws.Range("D45").Value = "10.64"
$ws.Range("D46").Value = "0.0533"
$ws.Range("D47").Value = "0.0227"
$ws.Range("D48").Value = "4.65"
$ws.Range("D49").Value = "1.931.72"
$ws.Range("D50").Value = "18.00"
$ws.Range("D51").Value = "115.10"

# Reset style back to default (no explicit style) now that text is stored,
# so the workbook formatting matches the original (unstyled) cells.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Coin name / link / volume columns (B, C, E) - plain text, safe to assign directly.
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +4.71%  "
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("E11").Value = "  +3.92%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("E26").Value = "  -3.34%  "
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("E35").Value = "  +5.03%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("E51").Value = "  +3.97%  "

Write-Output "Applied crypto list update"